$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    1 = 'basketball under pants'
    2 = 'softball gear for girls'
    3 = 'running capri'
    4 = 'softball compression sleeve'
    5 = 'youth softball compression sleeve'
    6 = 'running tights mens'
    7 = 'spandex men'
    8 = 'hockey kneepads'
    9 = 'padded leg sleeve'
    10 = 'mens basketball gear'
    11 = 'snowboarding padded shorts'
    12 = 'padded shorts snowboarding'
    13 = 'knee sleeve wrestling'
    14 = 'sleeve knee pads'
    15 = 'womens compression leggings'
    16 = 'airsoft knee pads'
    17 = 'mens compression tights 3 4'
    18 = 'basketball clothes for men'
    19 = 'men running tights'
    20 = 'knee pads nike'
    21 = 'knee pads mizuno'
    22 = 'knee pads bike'
    23 = 'yoga capri pants'
    24 = 'knee pads mtb'
    25 = 'knee pads skating'
    26 = 'mens workout tights'
    27 = 'mens basketball pants'
    28 = 'asics knee pads'
    29 = 'mens workout tights pants'
    30 = 'downhill knee pads'
    31 = 'men gym pants'
    32 = 'athletic capris'
    33 = 'valken knee pads'
    34 = 'woodland knee pads'
    35 = 'training tights men'
    36 = 'short tights for men'
    37 = 'ua compression pants'
    38 = 'men workout tights'
    39 = 'knee pads for exercise'
    40 = 'mens leggins'
    41 = 'nike kneepads'
    42 = 'youth football girdle with knee pads'
    43 = 'compression tights with pads'
    44 = 'compression leggings with knee pads'
    45 = 'basketball padded compression pants'
    46 = 'basketball knee pad pants'
    47 = 'padded compression pants men basketball'
    48 = 'padded tights men basketball'
    49 = 'tights with pads basketball'
    50 = 'mens leggings with knee pads'
    51 = 'leggings with knee pads women'
    52 = 'mtb knee pads men'
    53 = 'basketball padded knee sleeve'
    54 = 'compression pants women'
    55 = 'compression knee sleeve men basketball'
    56 = 'basketball sweat pants for men'
    57 = 'knee sleeve for wrestling'
    58 = 'leg sleeves for basketball youth'
    59 = 'training pants men'
    60 = 'compression knee sleeve men basketball'
    61 = 'basketball sweat pants for men'
    62 = 'knee sleeve for wrestling'
    63 = 'leg sleeves for basketball youth'
    64 = 'goalkeeper knee pads'
    65 = 'basketball calf sleeve'
    66 = 'compression knee sleeves with pads'
    67 = 'compression sleeve knee pads'
    68 = 'youth knee sleeve'
    69 = 'knee pad for scooter'
    70 = 'knee basketball'
    71 = 'knee pads for basketball youth'
    72 = 'mens compression knee'
    73 = 'knee pad sleeve basketball'
    74 = 'mens 3/4 compression pants'
    75 = 'youth compression knee pad sleeve'
    76 = 'mens basketball knee sleeves'
    77 = 'knee sleeve wrestling youth'
    78 = 'knee sleeves basketball men'
    79 = 'soccer compression pants'
    80 = 'leggings tight'
    81 = 'basketball leg sleeve youth padded'
    82 = 'knee pad construction'
    83 = 'youth basketball knee sleeve'
    84 = 'working knee pads for men'
    85 = 'cycling pants for men'
    86 = 'boys youth leggings'
    87 = 'compression running capris'
    88 = 'knee sleeve baseball'
    89 = 'compression knee sleeves for basketball'
    90 = 'volleyball kneepads'
    91 = 'compression knee sleeve with pad'
    92 = 'men capri shorts'
    93 = 'running compression pants'
    94 = 'mens work pants knee pads'
    95 = 'best knee pads'
    96 = 'compression pants sleeves'
    97 = 'mens compression running tights'
    98 = 'knee pads working'
    99 = 'basketball aids'
    100 = 'baseball youth pants'
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 1).Value = $values[$row]
}
